$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 473, shifting existing rows 473:495 down to 474:496.
$ws.Rows(473).Insert()

# Populate the newly-inserted row 473 with the new weekly record.
$ws.Cells.Item(473, 1).Value = 10
$ws.Cells.Item(473, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(473, 3).Value = "La Araucanía"
$ws.Cells.Item(473, 4).Value = 45041
$ws.Cells.Item(473, 5).Value = 9
$ws.Cells.Item(473, 6).Value = 100112009
$ws.Cells.Item(473, 7).Value = "Acelga"
$ws.Cells.Item(473, 8).Value = "Sin especificar"
$ws.Cells.Item(473, 9).Value = "Primera"
$ws.Cells.Item(473, 10).Value = 30
$ws.Cells.Item(473, 11).Value = 8000
$ws.Cells.Item(473, 12).Value = 8000
$ws.Cells.Item(473, 13).Value = 8000
$ws.Cells.Item(473, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(473, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(473, 16).Value = 667
$ws.Cells.Item(473, 17).Value = 12
$ws.Cells.Item(473, 18).Value = "Hortaliza"

# Match the date-style used by the rest of column D.
$ws.Cells.Item(473, 4).NumberFormat = $ws.Cells.Item(474, 4).NumberFormat
